$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting so that
# numeric-looking values (e.g. "1.002") are stored as text, not numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.808.52"
$ws.Range("E2").Value = "  +2.16%  "

$ws.Range("D3").Value = "2.118.36"
$ws.Range("E3").Value = "  +10.13%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "335.08"
$ws.Range("E5").Value = "  +4.82%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").Value = "0.5415"
$ws.Range("E7").Value = "  +6.63%  "

$ws.Range("D8").Value = "0.4424"
$ws.Range("E8").Value = "  +8.56%  "

$ws.Range("D9").Value = "0.09071"
$ws.Range("E9").Value = "  +8.64%  "

$ws.Range("D10").Value = "46.56"
$ws.Range("E10").Value = "  +9.88%  "

$ws.Range("D11").Value = "1.192"
$ws.Range("E11").Value = "  +6.92%  "

$ws.Range("D12").Value = "25.30"
$ws.Range("E12").Value = "  +5.25%  "

$ws.Range("D13").Value = "2.124.24"
$ws.Range("E13").Value = "  +10.43%  "

$ws.Range("D14").Value = "6.806"
$ws.Range("E14").Value = "  +5.66%  "

$ws.Range("D15").Value = "7.869"
$ws.Range("E15").Value = "  +8.43%  "

$ws.Range("D16").Value = "98.40"
$ws.Range("E16").Value = "  +6.18%  "

$ws.Range("D17").Value = "0.00001141"
$ws.Range("E17").Value = "  +4.13%  "

$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("D19").Value = "0.06652"
$ws.Range("E19").Value = "  +2.11%  "

$ws.Range("D20").Value = "19.33"
$ws.Range("E20").Value = "  +4.27%  "

$ws.Range("D21").Value = "6.430"
$ws.Range("E21").Value = "  +7.89%  "

$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("D23").Value = "30.911.53"
$ws.Range("E23").Value = "  +2.45%  "

$ws.Range("D24").Value = "12.17"
$ws.Range("E24").Value = "  +7.04%  "

$ws.Range("D25").Value = "2.375.41"
$ws.Range("E25").Value = "  +10.80%  "

$ws.Range("D26").Value = "2.275"
$ws.Range("E26").Value = "  +3.62%  "

$ws.Range("D27").Value = "22.97"
$ws.Range("E27").Value = "  +4.54%  "

$ws.Range("D28").Value = "2.579"
$ws.Range("E28").Value = "  +13.94%  "

$ws.Range("D29").Value = "163.71"
$ws.Range("E29").Value = "  +0.52%  "

$ws.Range("D30").Value = "134.39"
$ws.Range("E30").Value = "  +4.18%  "

$ws.Range("D31").Value = "1.177"
$ws.Range("E31").Value = "  +4.04%  "

$ws.Range("D32").Value = "0.1086"
$ws.Range("E32").Value = "  +3.77%  "

$ws.Range("D33").Value = "6.313"
$ws.Range("E33").Value = "  +5.79%  "

$ws.Range("D34").Value = "4.001"
$ws.Range("E34").Value = "  +5.26%  "

$ws.Range("D35").Value = "1.550"
$ws.Range("E35").Value = "  +27.65%  "

$ws.Range("D36").Value = "0.02602"
$ws.Range("E36").Value = "  +6.02%  "

$ws.Range("D37").Value = "5.602"
$ws.Range("E37").Value = "  +5.22%  "

$ws.Range("D38").Value = "9.633"
$ws.Range("E38").Value = "  +12.14%  "

$ws.Range("D39").Value = "0.06760"
$ws.Range("E39").Value = "  +4.96%  "

$ws.Range("D40").Value = "12.71"

$ws.Range("D41").Value = "0.2288"
$ws.Range("E41").Value = "  +6.37%  "

$ws.Range("D42").Value = "0.6867"
$ws.Range("E42").Value = "  +5.36%  "

$ws.Range("D43").Value = "1.261"
$ws.Range("E43").Value = "  +4.12%  "

$ws.Range("D44").Value = "14.19"
$ws.Range("E44").Value = "  +6.04%  "

$ws.Range("D45").Value = "0.6462"
$ws.Range("E45").Value = "  +6.76%  "

$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("D47").Value = "2.261"
$ws.Range("E47").Value = "  +3.58%  "

$ws.Range("D48").Value = "3.687"
$ws.Range("E48").Value = "  +1.59%  "

$ws.Range("D49").Value = "1.292"
$ws.Range("E49").Value = "  +6.60%  "

$ws.Range("D50").Value = "83.51"
$ws.Range("E50").Value = "  +6.18%  "

$ws.Range("D51").Value = "0.07097"
$ws.Range("E51").Value = "  +3.82%  "
